$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Populate the new "DT RMSE" column (D) with values for rows 2-31.
$dValues = @(
    20.05021,
    17.755682,
    9.2336460000000002,
    11.861848999999999,
    22.258886,
    17.437448,
    12.606989,
    12.095190000000001,
    15.726054,
    13.213855000000001,
    7.480315,
    19.219626999999999,
    4.8369559999999998,
    17.140916000000001,
    8.6910139999999991,
    9.1362299999999994,
    3.1621220000000001,
    7.1216160000000004,
    6.8786230000000002,
    8.1892469999999999,
    9.1766400000000008,
    10.32122,
    8.2251460000000005,
    5.7260390000000001,
    3.8187440000000001,
    3.5969799999999998,
    11.641442,
    5.2101040000000003,
    3.560991,
    12.730046
)

for ($i = 0; $i -lt $dValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}

# D32 already holds =AVERAGE(D2:D30); it will now recalculate to a real number
# now that D2:D30 contain data instead of being blank.

# 2. Apply a bordered-table look across the full A1:D32 range: an outer box
#    border plus medium dividers between every column.
$tableRange = $ws.Range("A1:D32")

# xlMedium = -4138, xlContinuous = 1
$edges = @(7, 8, 9, 10)  # xlEdgeLeft, xlEdgeTop, xlEdgeBottom, xlEdgeRight
foreach ($edge in $edges) {
    $tableRange.Borders.Item($edge).LineStyle = 1
    $tableRange.Borders.Item($edge).Weight = -4138
}

# xlInsideVertical = 11
$tableRange.Borders.Item(11).LineStyle = 1
$tableRange.Borders.Item(11).Weight = -4138

# 3. Update the selection stored in the worksheet view to the new table range.
$ws.Range("A1:D32").Select()
